$d = $word.ActiveDocument

function Get-ParagraphStartingWith($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($needle)) {
            return $p
        }
    }
    return $null
}

# --- 1) "3.1 Folder Lock (Availability)" -> strike through the whole paragraph ---
$p1 = Get-ParagraphStartingWith $d "3.1 Folder Lock"
$p1.Range.Font.StrikeThrough = 1

# --- 2) "3.6 Captcha (Co|nfidentiality)" / "3.7 Keylogger (Integrity)" ---
#        Re-join the two runs that made up "3.6 Captcha (Confidentiality)"
#        (currently split around the _GoBack bookmark) into a single run,
#        and relocate that _GoBack bookmark down into the "3.7 Keylogger"
#        paragraph, splitting it into "3.7" + " Keylogger (Integrity)"
#        around the bookmark instead. Splice both paragraphs' OOXML in one
#        go (over the combined range) so the bookmark can be moved exactly.
$p6 = Get-ParagraphStartingWith $d "3.6 Captcha"
$p7 = Get-ParagraphStartingWith $d "3.7 Keylogger"
$full = $d.Range($p6.Range.Start, $p7.Range.End)
$xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>3.6 Captcha (Confidentiality)</w:t></w:r></w:p><w:p><w:r><w:t>3.7</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> Keylogger (Integrity)</w:t></w:r></w:p></pkg:xmlData>'
[void]$full.InsertXML($xml)

# --- 3) "3.8 Website Blocker (Availability)" -> strike through the whole paragraph ---
$p8 = Get-ParagraphStartingWith $d "3.8 Website Blocker"
$p8.Range.Font.StrikeThrough = 1

Write-Host "Edit applied."
